$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "Coding Exercise 11: Checkerboard Exercise"
$ws.Range("E6").Formula = "=60+43"
$ws.Range("E16").Formula = "=SUM(E6:E8)"

$ws.Range("E7").Select()
